$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Prefix with an apostrophe so Excel stores these as literal text values
# (matching the original inline-string cells) instead of auto-converting
# numeric- or percent-looking text into numbers, which would also alter
# cell formatting/styles.

# Row 2 - AAPL
$ws.Range("B2").Value = "'173.80"
$ws.Range("C2").Value = "'+1.05"
$ws.Range("D2").Value = "'+0.61%"
$ws.Range("E2").Value = "'-1%"

# Row 3 - GOOG
$ws.Range("B3").Value = "'139.62"
$ws.Range("C3").Value = "'+0.68"
$ws.Range("D3").Value = "'+0.49%"
$ws.Range("E3").Value = "'6%"

# Row 4 - GOOGL
$ws.Range("B4").Value = "'138.50"
$ws.Range("C4").Value = "'+0.83"
$ws.Range("D4").Value = "'+0.60%"
$ws.Range("E4").Value = "'6%"

# Row 5 - AMZN
$ws.Range("B5").Value = "'175.39"
$ws.Range("C5").Value = "'+3.43"
$ws.Range("D5").Value = "'+1.99%"
$ws.Range("E5").Value = "'-5%"

# Row 6 - META
$ws.Range("B6").Value = "'483.59"
$ws.Range("C6").Value = "'-22.36"
$ws.Range("D6").Value = "'-4.42%"

# Row 7 - MSFT
$ws.Range("B7").Value = "'415.28"
$ws.Range("C7").Value = "'+10.76"
$ws.Range("D7").Value = "'+2.66%"
$ws.Range("E7").Value = "'-2%"

# Row 8 - NVDA
$ws.Range("B8").Value = "'919.13"
$ws.Range("C8").Value = "'+61.39"
$ws.Range("D8").Value = "'+7.16%"
$ws.Range("E8").Value = "'8%"
